$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values.
# A leading apostrophe forces Excel to keep the literal text, since many of
# these price strings (e.g. '50.00', '2.150') would otherwise be auto-converted
# to numbers and lose significant trailing zeros / formatting.

$ws.Range("D2").Value = "'23.449.99"
$ws.Range("E2").Value = "'  -1.40%  "
$ws.Range("D3").Value = "'1.646.38"
$ws.Range("E3").Value = "'  -0.53%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E5").Value = "'  +0.02%  "
$ws.Range("D6").Value = "'298.94"
$ws.Range("E6").Value = "'  -1.92%  "
$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = "'  -1.09%  "
$ws.Range("D8").Value = "'0.3534"
$ws.Range("E8").Value = "'  -2.21%  "
$ws.Range("D9").Value = "'50.00"
$ws.Range("E9").Value = "'  -2.45%  "
$ws.Range("D10").Value = "'0.08082"
$ws.Range("E10").Value = "'  -1.74%  "
$ws.Range("E11").Value = "'  -3.10%  "
$ws.Range("E12").Value = "'  -0.04%  "
$ws.Range("D13").Value = "'22.10"
$ws.Range("E13").Value = "'  -2.70%  "
$ws.Range("D14").Value = "'6.391"
$ws.Range("E14").Value = "'  -2.32%  "
$ws.Range("D15").Value = "'7.331"
$ws.Range("E15").Value = "'  -1.56%  "
$ws.Range("E16").Value = "'  -3.13%  "
$ws.Range("D17").Value = "'1.650.28"
$ws.Range("E17").Value = "'  +0.55%  "
$ws.Range("D18").Value = "'96.96"
$ws.Range("E18").Value = "'  -0.89%  "
$ws.Range("D19").Value = "'0.06962"
$ws.Range("E19").Value = "'  -0.34%  "
$ws.Range("D20").Value = "'6.752"
$ws.Range("E20").Value = "'  -0.36%  "
$ws.Range("D21").Value = "'17.41"
$ws.Range("E21").Value = "'  -2.07%  "
$ws.Range("E22").Value = "'  -0.04%  "
$ws.Range("D23").Value = "'12.43"
$ws.Range("E23").Value = "'  -2.41%  "
$ws.Range("D24").Value = "'23.477.04"
$ws.Range("E24").Value = "'  -1.30%  "
$ws.Range("D25").Value = "'2.499"
$ws.Range("E25").Value = "'  -1.82%  "
$ws.Range("D26").Value = "'2.886"
$ws.Range("E26").Value = "'  -6.68%  "
$ws.Range("D27").Value = "'20.88"
$ws.Range("E27").Value = "'  -2.16%  "
$ws.Range("D28").Value = "'152.26"
$ws.Range("D29").Value = "'5.192"
$ws.Range("E29").Value = "'  -0.97%  "
$ws.Range("D30").Value = "'132.75"
$ws.Range("E30").Value = "'  -1.73%  "
$ws.Range("D31").Value = "'1.830.01"
$ws.Range("E31").Value = "'  +0.16%  "
$ws.Range("D32").Value = "'6.955"
$ws.Range("E32").Value = "'  +0.73%  "
$ws.Range("D33").Value = "'2.150"
$ws.Range("E33").Value = "'  +1.00%  "
$ws.Range("D34").Value = "'11.43"
$ws.Range("E34").Value = "'  -4.02%  "
$ws.Range("D35").Value = "'0.9910"
$ws.Range("E35").Value = "'  -9.23%  "
$ws.Range("D36").Value = "'0.02714"
$ws.Range("E36").Value = "'  -4.69%  "
$ws.Range("D37").Value = "'0.08720"
$ws.Range("E37").Value = "'  -1.42%  "
$ws.Range("E38").Value = "'  -3.38%  "
$ws.Range("D39").Value = "'5.936"
$ws.Range("E39").Value = "'  -3.48%  "
$ws.Range("D40").Value = "'12.92"
$ws.Range("E40").Value = "'  -0.36%  "
$ws.Range("D41").Value = "'0.06793"
$ws.Range("E41").Value = "'  -6.01%  "
$ws.Range("D42").Value = "'0.6883"
$ws.Range("E42").Value = "'  -2.90%  "
$ws.Range("D43").Value = "'1.294"
$ws.Range("E43").Value = "'  -3.61%  "
$ws.Range("D44").Value = "'15.68"
$ws.Range("E44").Value = "'  -1.29%  "
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("D46").Value = "'0.6357"
$ws.Range("E46").Value = "'  -3.09%  "
$ws.Range("D47").Value = "'2.253"
$ws.Range("E47").Value = "'  -3.65%  "
$ws.Range("D48").Value = "'3.906"
$ws.Range("E48").Value = "'  -1.63%  "
$ws.Range("D49").Value = "'0.07723"
$ws.Range("E49").Value = "'  -3.32%  "
$ws.Range("D50").Value = "'127.51"
$ws.Range("E50").Value = "'  -1.02%  "
$ws.Range("D51").Value = "'1.150"
$ws.Range("E51").Value = "'  -4.03%  "
